$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New test case: LoginWithValidPasswordAsEmail (row 7) ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "LoginWithValidPasswordAsEmail"
$ws.Range("C7").Value = "Existing"

# Write B9's text now so the shared-string table gets the new unique
# strings in the same order the original workbook has them in
# (LoginWithRandomNumbersAsEmailAndPassword before LoginWithInvalidEmail).
$ws.Range("B9").Value = "LoginWithRandomNumbersAsEmailAndPassword"

# --- New test case: LoginWithInvalidEmail (row 8) ---
# Email column keeps a mailto: hyperlink (like the other test rows) that
# displays/links to abv@abv.bg, but the underlying cell value is the
# numeric "invalid email" 1234.
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "LoginWithInvalidEmail"
$ws.Hyperlinks.Add($ws.Range("C8"), "mailto:abv@abv.bg", "", "", "abv@abv.bg")
$ws.Range("C8").Value = 1234
$ws.Range("C8").Style = "Hyperlink"

# --- New test case: LoginWithRandomNumbersAsEmailAndPassword (row 9) ---
$ws.Range("A9").Value = 8
$ws.Range("C9").Value = 1234
$ws.Range("D9").Value = 1234

# Matches the author's final selection being on the last edited cell.
$ws.Range("D9").Select()
